$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four more student-id / is_Delete rows below the existing data
$ws.Range("A3").Value = 3121410488
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 3121410458
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = 3121410459
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = 3121410460
$ws.Range("B6").Value = 0

# Move the active selection to the next empty row, like the saved file shows
$ws.Range("B7").Select() | Out-Null
